$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source "Hortaliza, Vega Modelo de Temuco - Albahaca" sheet gained one new
# weekly price record. It was inserted as row 368, pushing the previously
# existing rows 368-409 down to 369-410 (dimension grows from A1:R409 to
# A1:R410).
$ws.Rows(368).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A368").Value2 = 10
$ws.Range("B368").Value2 = "Vega Modelo de Temuco"
$ws.Range("C368").Value2 = "La Araucanía"
$ws.Range("D368").Value2 = 45124
$ws.Range("E368").Value2 = 9
$ws.Range("F368").Value2 = 100112052
$ws.Range("G368").Value2 = "Albahaca"
$ws.Range("H368").Value2 = "Sin especificar"
$ws.Range("I368").Value2 = "Primera"
$ws.Range("J368").Value2 = 250
$ws.Range("K368").Value2 = 5000
$ws.Range("L368").Value2 = 6000
$ws.Range("M368").Value2 = 5400
$ws.Range("N368").Value2 = "$/paquete"
$ws.Range("O368").Value2 = "Región de Arica y Parinacota"
$ws.Range("P368").Value2 = 5400
$ws.Range("Q368").Value2 = 1
$ws.Range("R368").Value2 = "Hortaliza"
